# Updated cryptos list (price/volume refresh).
# Note: some Price (column D) values look like plain numbers (e.g. "136.12").
# These cells must stay TEXT (as in the original sheet), so such values are
# written with a leading apostrophe ('') to force Excel to keep them as text
# instead of silently converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.866.07'
$ws.Range('E2').Value = '  -3.33%  '
$ws.Range('D3').Value = '3.229.85'
$ws.Range('E3').Value = '  -3.90%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''540.07'
$ws.Range('E5').Value = '  -4.72%  '
$ws.Range('D6').Value = '''136.12'
$ws.Range('E6').Value = '  -8.40%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '3.228.53'
$ws.Range('E8').Value = '  -3.91%  '
$ws.Range('E9').Value = '  -4.31%  '
$ws.Range('D10').Value = '''7.61'
$ws.Range('E10').Value = '  -4.13%  '
$ws.Range('E11').Value = '  -5.77%  '
$ws.Range('D12').Value = '''0.394'
$ws.Range('E12').Value = '  -4.44%  '
$ws.Range('D13').Value = '3.784.72'
$ws.Range('E13').Value = '  -3.80%  '
$ws.Range('E14').Value = '  -0.94%  '
$ws.Range('D15').Value = '''25.99'
$ws.Range('E15').Value = '  -7.06%  '
$ws.Range('D16').Value = '3.235.60'
$ws.Range('E16').Value = '  -3.65%  '
$ws.Range('E17').Value = '  -5.96%  '
$ws.Range('D18').Value = '58.901.29'
$ws.Range('E18').Value = '  -3.45%  '
$ws.Range('D19').Value = '''5.89'
$ws.Range('E19').Value = '  -7.20%  '
$ws.Range('D20').Value = '''13.36'
$ws.Range('E20').Value = '  -5.82%  '
$ws.Range('D21').Value = '''8.25'
$ws.Range('E21').Value = '  -6.48%  '
$ws.Range('D22').Value = '''362.25'
$ws.Range('E22').Value = '  -3.17%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').Value = '''70.55'
$ws.Range('E24').Value = '  -6.26%  '
$ws.Range('D25').Value = '''0.520'
$ws.Range('E25').Value = '  -6.86%  '
$ws.Range('D26').Value = '3.368.32'
$ws.Range('E26').Value = '  -3.75%  '
$ws.Range('D27').Value = '''0.170'
$ws.Range('E27').Value = '  -2.77%  '
$ws.Range('D28').Value = '0.0₃0968'
$ws.Range('E28').Value = '  -10.87%  '
$ws.Range('D29').Value = '''0.998'
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('D30').Value = '''7.08'
$ws.Range('E30').Value = '  -4.01%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').Value = '''1.93'
$ws.Range('E32').Value = '  -6.73%  '
$ws.Range('D33').Value = '''7.10'
$ws.Range('E33').Value = '  -7.64%  '
$ws.Range('D34').Value = '''21.93'
$ws.Range('E34').Value = '  -4.08%  '
$ws.Range('D35').Value = '''1.22'
$ws.Range('E35').Value = '  -5.47%  '
$ws.Range('D36').Value = '''4.95'
$ws.Range('E36').Value = '  -7.95%  '
$ws.Range('D37').Value = '''162.88'
$ws.Range('E37').Value = '  -3.47%  '
$ws.Range('D38').Value = '''6.42'
$ws.Range('E38').Value = '  -5.30%  '
$ws.Range('E39').Value = '  -6.77%  '
$ws.Range('D40').Value = '''26.27'
$ws.Range('E40').Value = '  -9.76%  '
$ws.Range('D41').Value = '''0.0709'
$ws.Range('E41').Value = '  -4.90%  '
$ws.Range('D42').Value = '3.261.85'
$ws.Range('E42').Value = '  -3.91%  '
$ws.Range('D43').Value = '''41.05'
$ws.Range('E43').Value = '  -2.94%  '
$ws.Range('D44').Value = '''0.717'
$ws.Range('E44').Value = '  -5.61%  '
$ws.Range('D45').Value = '''1.10'
$ws.Range('E45').Value = '  -3.47%  '
$ws.Range('D46').Value = '''4.03'
$ws.Range('E46').Value = '  -5.90%  '
$ws.Range('E47').Value = '  -6.27%  '
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').Value = '2.299.94'
$ws.Range('E49').Value = '  -7.69%  '
$ws.Range('E50').Value = '  -5.30%  '
$ws.Range('D51').Value = '''20.94'
$ws.Range('E51').Value = '  -7.17%  '
